$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 949, shifting existing rows 949:1067 down to 950:1068
$ws.Rows("949:949").Insert()

# Populate the newly inserted row 949 with the new record data
$ws.Range("A949").Value = 3
$ws.Range("B949").Value = "Femacal de La Calera"
$ws.Range("C949").Value = "Coquimbo"
$ws.Range("D949").Value = 45142
$ws.Range("E949").Value = 5
$ws.Range("F949").Value = 100112006
$ws.Range("G949").Value = "Repollo"
$ws.Range("H949").Value = "Crespo record"
$ws.Range("I949").Value = "Primera"
$ws.Range("J949").Value = 2400
$ws.Range("K949").Value = 750
$ws.Range("L949").Value = 800
$ws.Range("M949").Value = 775
$ws.Range("N949").Value = "`$/unidad"
$ws.Range("O949").Value = "Provincia de Quillota"
$ws.Range("P949").Value = 775
$ws.Range("Q949").Value = 1
$ws.Range("R949").Value = "Hortaliza"
